$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Significance threshold" text box and update its wording to
# "Significant transcripts" (David's review comment for lecture 17).
$updated = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "Significance threshold") {
            $tr.Text = "Significant transcripts"
            $updated = $true
        }
    }
}

if (-not $updated) {
    throw "Could not find the 'Significance threshold' text box to update"
}
